$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 22289
$ws.Range("E2").Value = 478
$ws.Range("F2").Value = 478
$ws.Range("G2").Value = 520
$ws.Range("H2").Value = 437
$ws.Range("I2").Value = 437
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 182217
$ws.Range("L2").Value = 165788
$ws.Range("M2").Value = 16429
$ws.Range("N2").Value = 16428
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 4349
$ws.Range("Q2").Value = -2265
$ws.Range("R2").Value = -3033
$ws.Range("S2").Value = 5096
$ws.Range("T2").Value = 971
$ws.Range("V2").Value = 28657
$ws.Range("W2").Value = 2.14
$ws.Range("X2").Value = 1.96
$ws.Range("Y2").Value = 2.67
$ws.Range("Z2").Value = 0.27
$ws.Range("AA2").Value = 1009.12
$ws.Range("AB2").Value = 315.53
$ws.Range("AC2").Value = 503
$ws.Range("AD2").Value = 18.85
$ws.Range("AE2").Value = 22404
$ws.Range("AF2").Value = 0.42
$ws.Range("AG2").Value = 250
$ws.Range("AH2").Value = 2.63
$ws.Range("AI2").Value = 44.48
$ws.Range("AJ2").Value = 50773400
$ws.Range("U2").ClearContents()

# Row 3
$ws.Range("D3").Value = 34411
$ws.Range("E3").Value = 1701
$ws.Range("F3").Value = 1701
$ws.Range("G3").Value = 1807
$ws.Range("H3").Value = 1362
$ws.Range("I3").Value = 1362
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 186860
$ws.Range("L3").Value = 169545
$ws.Range("M3").Value = 17315
$ws.Range("N3").Value = 17314
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 4349
$ws.Range("Q3").Value = 3424
$ws.Range("R3").Value = -347
$ws.Range("S3").Value = -3092
$ws.Range("T3").Value = 488
$ws.Range("V3").Value = 29671
$ws.Range("W3").Value = 4.94
$ws.Range("X3").Value = 3.96
$ws.Range("Y3").Value = 8.07
$ws.Range("Z3").Value = 0.74
$ws.Range("AA3").Value = 979.17
$ws.Range("AB3").Value = 340.94
$ws.Range("AC3").Value = 1569
$ws.Range("AD3").Value = 7.55
$ws.Range("AE3").Value = 24390
$ws.Range("AF3").Value = 0.49
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 4.22
$ws.Range("AI3").Value = 26.85
$ws.Range("AJ3").Value = 50773400
$ws.Range("U3").ClearContents()

# Row 4
$ws.Range("D4").Value = 41346
$ws.Range("E4").Value = 833
$ws.Range("F4").Value = 833
$ws.Range("G4").Value = 861
$ws.Range("H4").Value = 740
$ws.Range("I4").Value = 740
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 176926
$ws.Range("L4").Value = 159145
$ws.Range("M4").Value = 17781
$ws.Range("N4").Value = 17780
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 4349
$ws.Range("Q4").Value = -7944
$ws.Range("R4").Value = 1953
$ws.Range("S4").Value = 7190
$ws.Range("T4").Value = 1013
$ws.Range("V4").Value = 33973
$ws.Range("W4").Value = 2.02
$ws.Range("X4").Value = 1.79
$ws.Range("Y4").Value = 4.22
$ws.Range("Z4").Value = 0.41
$ws.Range("AA4").Value = 895.05
$ws.Range("AB4").Value = 350.75
$ws.Range("AC4").Value = 853
$ws.Range("AD4").Value = 12.19
$ws.Range("AE4").Value = 24952
$ws.Range("AF4").Value = 0.42
$ws.Range("AG4").Value = 550
$ws.Range("AH4").Value = 5.29
$ws.Range("AI4").Value = 54.39
$ws.Range("AJ4").Value = 50773400
$ws.Range("U4").ClearContents()

# Row 5
$ws.Range("D5").Value = 35484
$ws.Range("E5").Value = 1330
$ws.Range("F5").Value = 1330
$ws.Range("G5").Value = 1513
$ws.Range("H5").Value = 1159
$ws.Range("I5").Value = 1159
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 174468
$ws.Range("L5").Value = 155789
$ws.Range("M5").Value = 18679
$ws.Range("N5").Value = 18677
$ws.Range("O5").Value = 2
$ws.Range("P5").Value = 4349
$ws.Range("Q5").Value = -6436
$ws.Range("R5").Value = -1866
$ws.Range("S5").Value = 8439
$ws.Range("T5").Value = 767
$ws.Range("V5").Value = 41860
$ws.Range("W5").Value = 3.75
$ws.Range("X5").Value = 3.27
$ws.Range("Y5").Value = 6.36
$ws.Range("Z5").Value = 0.66
$ws.Range("AA5").Value = 834.04
$ws.Range("AB5").Value = 370.72
$ws.Range("AC5").Value = 1336
$ws.Range("AD5").Value = 10.78
$ws.Range("AE5").Value = 26135
$ws.Range("AF5").Value = 0.55
$ws.Range("AG5").Value = 610
$ws.Range("AH5").Value = 4.24
$ws.Range("AI5").Value = 38.54
$ws.Range("AJ5").Value = 50773400
$ws.Range("U5").ClearContents()

# Row 6
$ws.Range("D6").Value = 25570
$ws.Range("E6").Value = 1585
$ws.Range("F6").Value = 1585
$ws.Range("G6").Value = 1893
$ws.Range("H6").Value = 1407
$ws.Range("I6").Value = 1408
$ws.Range("K6").Value = 191513
$ws.Range("L6").Value = 171504
$ws.Range("M6").Value = 20009
$ws.Range("N6").Value = 19692
$ws.Range("P6").Value = 4349
$ws.Range("Q6").Value = -15981
$ws.Range("R6").Value = -98
$ws.Range("S6").Value = 18935
$ws.Range("T6").Value = 191
$ws.Range("V6").Value = 54370
$ws.Range("W6").Value = 6.2
$ws.Range("X6").Value = 5.5
$ws.Range("Y6").Value = 7.34
$ws.Range("Z6").Value = 0.77
$ws.Range("AA6").Value = 857.14
$ws.Range("AB6").Value = 400.48
$ws.Range("AC6").Value = 1623
$ws.Range("AD6").Value = 6.99
$ws.Range("AE6").Value = 27460
$ws.Range("AF6").Value = 0.41
$ws.Range("AG6").Value = 620
$ws.Range("AH6").Value = 5.46
$ws.Range("AI6").Value = 32.33
$ws.Range("AJ6").Value = 50773400
$ws.Range("U6").ClearContents()

# Row 7
$ws.Range("E7").Value = 1120
$ws.Range("G7").Value = 1440
$ws.Range("H7").Value = 1090
$ws.Range("I7").Value = 1080
$ws.Range("K7").Value = 205810
$ws.Range("L7").Value = 185190
$ws.Range("M7").Value = 20620
$ws.Range("N7").Value = 20290
$ws.Range("P7").Value = 4350
$ws.Range("Y7").Value = 5.4
$ws.Range("Z7").Value = 0.55
$ws.Range("AA7").Value = 898.11
$ws.Range("AC7").Value = 1245
$ws.Range("AD7").Value = 8.56
$ws.Range("AE7").Value = 29852
$ws.Range("AF7").Value = 0.36
$ws.Range("AG7").Value = 620
$ws.Range("AH7").Value = 5.82
$ws.Range("AI7").Value = 29.15
$ws.Range("D7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()

# Row 8
$ws.Range("E8").Value = 1380
$ws.Range("G8").Value = 1350
$ws.Range("H8").Value = 1020
$ws.Range("I8").Value = 1010
$ws.Range("K8").Value = 209090
$ws.Range("L8").Value = 187900
$ws.Range("M8").Value = 21190
$ws.Range("N8").Value = 20860
$ws.Range("P8").Value = 4350
$ws.Range("Y8").Value = 4.91
$ws.Range("Z8").Value = 0.49
$ws.Range("AA8").Value = 886.74
$ws.Range("AC8").Value = 1164
$ws.Range("AD8").Value = 9.15
$ws.Range("AE8").Value = 30690
$ws.Range("AF8").Value = 0.35
$ws.Range("AG8").Value = 620
$ws.Range("AH8").Value = 5.82
$ws.Range("AI8").Value = 31.17
$ws.Range("D8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()

# Row 9
$ws.Range("E9").Value = 1480
$ws.Range("G9").Value = 1450
$ws.Range("H9").Value = 1100
$ws.Range("I9").Value = 1090
$ws.Range("K9").Value = 208130
$ws.Range("L9").Value = 186290
$ws.Range("M9").Value = 21840
$ws.Range("N9").Value = 21500
$ws.Range("P9").Value = 4350
$ws.Range("Y9").Value = 5.15
$ws.Range("Z9").Value = 0.53
$ws.Range("AA9").Value = 852.98
$ws.Range("AC9").Value = 1256
$ws.Range("AD9").Value = 8.48
$ws.Range("AE9").Value = 31632
$ws.Range("AF9").Value = 0.34
$ws.Range("AG9").Value = 630
$ws.Range("AH9").Value = 5.92
$ws.Range("AI9").Value = 29.35
$ws.Range("D9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
